$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing rows (40-64) that no longer exist in the updated data
$ws.Range("A40:A64").EntireRow.Delete()

# Update cluster names and active-case counts for the remaining rows (2-39)
$ws.Range("A2").Value = "3323 Villa Maria Catholic Homes St Bernadette's Aged Care Sunshine North"
$ws.Range("B2").Value = 10
$ws.Range("A3").Value = "3376 Royal Freemasons Coppin Centre Melbourne"
$ws.Range("B3").Value = 11
$ws.Range("A4").Value = "3601 Baptcare Westhaven community"
$ws.Range("B4").Value = 13
$ws.Range("A5").Value = "3653 Fronditha Thalpori St Albans Aged Care"
$ws.Range("B5").Value = 38
$ws.Range("A6").Value = "44121 Wallan Primary School Wallan Outbreak"
$ws.Range("B6").Value = 11
$ws.Range("A7").Value = "44165 Greenvale Primary School"
$ws.Range("B7").Value = 11
$ws.Range("A8").Value = "44226 Boneo Primary School Boneo"
$ws.Range("B8").Value = 10
$ws.Range("A9").Value = "44234 Lucknow Primary School"
$ws.Range("B9").Value = 12
$ws.Range("A10").Value = "44321 Maiden Gully Primary School Maiden Gully"
$ws.Range("B10").Value = 11
$ws.Range("A11").Value = "44395 Buln Buln Primary School"
$ws.Range("B11").Value = 10
$ws.Range("A12").Value = "44852 Dandenong South Primary School Dandenong"
$ws.Range("B12").Value = 14
$ws.Range("A13").Value = "44978 Deer Park West Primary School Deer Park"
$ws.Range("B13").Value = 10
$ws.Range("A14").Value = "45034 River Gum Primary School Hampton Park"
$ws.Range("B14").Value = 11
$ws.Range("A15").Value = "45158 Rowellyn Park Primary School Carrum Downs"
$ws.Range("B15").Value = 14
$ws.Range("A16").Value = "45161 Woodlands Primary School Langwarrin"
$ws.Range("B16").Value = 10
$ws.Range("A17").Value = "45249 Creekside K-9 College Caroline Springs"
$ws.Range("B17").Value = 13
$ws.Range("A18").Value = "45585 Mount Ridley College Craigieburn"
$ws.Range("B18").Value = 14
$ws.Range("A19").Value = "45695 Sacred Heart Primary School Yarrawonga Outbreak"
$ws.Range("B19").Value = 54
$ws.Range("A20").Value = "4574 Village Glen Aged Care Residences Mornington"
$ws.Range("B20").Value = 16
$ws.Range("A21").Value = "45804 St Therese's School Essendon"
$ws.Range("B21").Value = 13
$ws.Range("A22").Value = "45809 St Finbar's Primary School Brighton East"
$ws.Range("B22").Value = 12
$ws.Range("A23").Value = "46050 Our Lady's Catholic Primary School Craigieburn"
$ws.Range("B23").Value = 26
$ws.Range("A24").Value = "46322 Minaret College Officer Campus Officer"
$ws.Range("B24").Value = 20
$ws.Range("A25").Value = "46328 Ilim College Dallas Primary Campus Tier 1A Dallas"
$ws.Range("B25").Value = 11
$ws.Range("A26").Value = "46390 Al Siraat College Epping"
$ws.Range("B26").Value = 44
$ws.Range("A27").Value = "51478 Wyndham Vale Primary School Wyndham Vale"
$ws.Range("B27").Value = 10
$ws.Range("A28").Value = "52380 Al Iman College Melton South"
$ws.Range("B28").Value = 28
$ws.Range("A29").Value = "52786 Hume Anglican Grammar Donnybrook Campus"
$ws.Range("B29").Value = 16
$ws.Range("A30").Value = "Adass Israel School Elsternwick"
$ws.Range("B30").Value = 11
$ws.Range("A31").Value = "Ilim College Dallas Main Campus Dallas Oct"
$ws.Range("B31").Value = 24
$ws.Range("A32").Value = "Ilim College Kiewa Campus Dallas"
$ws.Range("B32").Value = 12
$ws.Range("A33").Value = "Inverloch Primary School"
$ws.Range("B33").Value = 12
$ws.Range("A34").Value = "Islamic College of Melbourne Tarneit Oct Nov"
$ws.Range("B34").Value = 26
$ws.Range("A35").Value = "Lavalla Catholic College St Pauls Campus Traralgon"
$ws.Range("B35").Value = 11
$ws.Range("A36").Value = "Middle Park Primary School Middle Park"
$ws.Range("B36").Value = 12
$ws.Range("A37").Value = "Morwell Park Primary School Morwell"
$ws.Range("B37").Value = 29
$ws.Range("A38").Value = "Nio Early Learning Adventures Preston"
$ws.Range("B38").Value = 20
$ws.Range("A39").Value = "Supreme Caravans Manufacturing Campbellfield"
$ws.Range("B39").Value = 37
